$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (row, newValue) for column F ("想去人数")
$updates = @{
    "展览"     = @{ 2 = 53; 3 = 21532; 8 = 8011; 19 = 1364; 20 = 554; 22 = 712; 34 = 5124; 39 = 13241 }
    "全部类型" = @{ 2 = 53; 3 = 21532; 6 = 8011; 16 = 1364; 17 = 554; 19 = 712; 34 = 5124; 39 = 13241 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $newValue = $rowsMap[$row]
        $ws.Range("F$row").Value = $newValue
    }
}

$wb.Save()
